$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 2.88
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 2.6
$ws.Range("I2").Value = 2.64
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 3.75
$ws.Range("P2").Value = 2.14
$ws.Range("Q2").Value = 1.78
$ws.Range("S2").Value = 2.96
$ws.Range("T2").Value = 1.66
$ws.Range("V2").Value = 1.61
$ws.Range("W2").Value = 1.51
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 21
$ws.Range("AD2").Value = 12
$ws.Range("AF2").Value = 23
$ws.Range("AG2").Value = 14
$ws.Range("AH2").Value = 16.5
$ws.Range("AI2").Value = 36
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 24

# Row 4
$ws.Range("F4").Value = 2.62
$ws.Range("G4").Value = 3.25
$ws.Range("H4").Value = 2.46
$ws.Range("I4").Value = 2.86
$ws.Range("K4").Value = 4.1

# Row 5
$ws.Range("G5").Value = 4.7

# Row 6
$ws.Range("H6").Value = 2.18
$ws.Range("I6").Value = 2.2
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 3.65
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 1.96
$ws.Range("T6").Value = 1.76
$ws.Range("U6").Value = 2.18
$ws.Range("X6").Value = 15.5
$ws.Range("Y6").Value = 10
$ws.Range("AA6").Value = 32
$ws.Range("AB6").Value = 14.5
$ws.Range("AC6").Value = 8.4
$ws.Range("AE6").Value = 25
$ws.Range("AF6").Value = 32
$ws.Range("AG6").Value = 15.5
$ws.Range("AH6").Value = 17.5
$ws.Range("AI6").Value = 36
$ws.Range("AJ6").Value = 85
$ws.Range("AK6").Value = 55
$ws.Range("AO6").Value = 16.5

# Row 7
$ws.Range("F7").Value = 2.24
$ws.Range("G7").Value = 2.68
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 3.85
$ws.Range("J7").Value = 3.1
$ws.Range("K7").Value = 3.9
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.89

# Row 8
$ws.Range("G8").Value = 2.68
$ws.Range("H8").Value = 2.52

# Row 10
$ws.Range("G10").Value = 3.25

# Row 12
$ws.Range("G12").Value = 2.1
$ws.Range("H12").Value = 3.75
$ws.Range("J12").Value = 3.9
$ws.Range("K12").Value = 4
$ws.Range("P12").Value = 2.3
$ws.Range("T12").Value = 1.66
$ws.Range("X12").Value = 20
$ws.Range("Y12").Value = 18.5
$ws.Range("Z12").Value = 32
$ws.Range("AF12").Value = 14.5
$ws.Range("AJ12").Value = 26
$ws.Range("AL12").Value = 32
$ws.Range("AN12").Value = 11.5
$ws.Range("AO12").Value = 38

# Row 13
$ws.Range("F13").Value = 6.2
$ws.Range("G13").Value = 6.4
$ws.Range("H13").Value = 1.57
$ws.Range("I13").Value = 1.59
$ws.Range("K13").Value = 5
$ws.Range("P13").Value = 2.66
$ws.Range("T13").Value = 1.69
$ws.Range("U13").Value = 2.34
$ws.Range("AA13").Value = 15.5
$ws.Range("AB13").Value = 30
$ws.Range("AE13").Value = 15
$ws.Range("AG13").Value = 25
$ws.Range("AJ13").Value = 160
$ws.Range("AK13").Value = 75
$ws.Range("AL13").Value = 65
$ws.Range("AN13").Value = 65

# Row 14
$ws.Range("F14").Value = 3.8
$ws.Range("G14").Value = 3.9
$ws.Range("H14").Value = 1.97
$ws.Range("I14").Value = 1.99
$ws.Range("P14").Value = 2.6
$ws.Range("U14").Value = 2.6
$ws.Range("Y14").Value = 13.5
$ws.Range("AA14").Value = 25
$ws.Range("AF14").Value = 36
$ws.Range("AG14").Value = 16.5
$ws.Range("AH14").Value = 16
$ws.Range("AJ14").Value = 70
$ws.Range("AK14").Value = 40
$ws.Range("AM14").Value = 60
$ws.Range("AN14").Value = 28

# Row 15
$ws.Range("F15").Value = 1.9
$ws.Range("G15").Value = 2.02
$ws.Range("H15").Value = 4.8
$ws.Range("I15").Value = 5.4
$ws.Range("K15").Value = 3.55
$ws.Range("Q15").Value = 2.44
